$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New diary entry row (row 37), scaffolded the same way as the existing rows.
$row = 37

$ws.Cells.Item($row, 1).Value = "28 marras"
$ws.Cells.Item($row, 2).Value = "12.00-13.30, 19.15-20.15"
$ws.Cells.Item($row, 3).Value = "Tsemppirinki, Millingtonin törmäysdemon pohjakoodiin tutustumista"
$ws.Cells.Item($row, 7).Value = 2.5

# Match formatting of the row above it (wrap text on B/C columns, row height).
$ws.Cells.Item($row, 2).NumberFormat = "h:mm"
$ws.Cells.Item($row, 2).WrapText = $true
$ws.Cells.Item($row, 3).WrapText = $true
$ws.Rows.Item($row).RowHeight = 43.2

$ws.Range("A39").Select()
